$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows at row 7 (old row 7 -> 10, old row 10 -> 13, etc.)
$ws.Rows("7:9").Insert()

# --- Row 4: C4 picks up the same blank/formatted style as F3 ---
$ws.Range("F3").Copy()
$ws.Range("C4").PasteSpecial(-4122)   # xlPasteFormats

# --- Row 6: relabel + new formula (divide Neg thickness by electrode count) ---
$ws.Range("A6").Value = "Electrolyte Neg Step Size"
$ws.Range("B6").Formula = "=A3/B9"

# NB: shared-string table order matters for the diff, so the "Pos" label is
# entered before the "Sep" label (matches the author's original edit order).
# --- Row 8: "Electrolyte Pos Step Size" ---
$ws.Range("A8").Value = "Electrolyte Pos Step Size"
$ws.Range("B8").Formula = "=C3/B9"

# --- Row 7: "Electrolyte Sep Step Size" ---
$ws.Range("A7").Value = "Electrolyte Sep Step Size"
$ws.Range("B7").Formula = "=B3/B10"

# --- Row 9: "Number in Electrodes" (plain number, default/general format) ---
$ws.Range("A9").Value = "Number in Electrodes"
$ws.Range("B9").Value = 3
$ws.Range("B9").ClearFormats()

# --- Row 10: "Number in Seperator" (plain number, default/general format) ---
$ws.Range("A10").Value = "Number in Seperator"
$ws.Range("B10").Value = 1
$ws.Range("B10").ClearFormats()

# --- Row 12: new blank, formatted cells above the step table ---
$ws.Range("E12:H12").NumberFormat = "0.000000E+00"

# --- Row 13 (former step-table row, now fed by the 3 separate step sizes) ---
$ws.Range("D13").Formula = "=C13+B6"
$ws.Range("E13").Formula = "=D13+B7"
$ws.Range("F13").Formula = "=B8+E13"
$ws.Range("G13").Formula = "=F13+B8"
$ws.Range("H13").Formula = "=G13+B8"
$ws.Range("B13:C13").Formula = "=A13+`$B`$6"
$ws.Range("A13:H13").NumberFormat = "0.000E+00"

$excel.Calculate()

# --- Column widths: A widened for the longer labels; D gets its own bestFit ---
$ws.Range("A1").ColumnWidth = 22.0
$ws.Range("D1").ColumnWidth = 11.833333333333334

# --- Active cell / selection moves to E10 ---
$ws.Range("E10").Select()

Write-Output "stage1 done"
